$d = $word.ActiveDocument

# The <id>...</id> markup around each p157r_N label used to be split
# across three separate runs ("<id>", "p157r_N", "</id>"). Collapse each
# one back into a single run ("<id>p157r_N</id>") by finding the whole
# visible string and replacing it with itself - Word merges the runs
# that the found range spans into one run using the first run's
# formatting.
for ($i = 1; $i -le 5; $i++) {
    $needle = "<id>p157r_$i</id>"
    $rng = $d.Content
    $found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, `
                                $true, 1, $false, $needle, 2)
}
